# Auto-generated edit script: updates cryptos list values to match the commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.183.71"
$ws.Range("E2").Value = "  -2.10%  "

$ws.Range("D3").Value = "2.581.91"
$ws.Range("E3").Value = "  -2.45%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "'563.64"
$ws.Range("E5").Value = "  -0.96%  "

$ws.Range("D6").Value = "'142.48"
$ws.Range("E6").Value = "  -2.09%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "'0.596"
$ws.Range("E8").Value = "  -1.81%  "

$ws.Range("D9").Value = "2.588.41"
$ws.Range("E9").Value = "  -2.70%  "

$ws.Range("D10").Value = "'6.63"
$ws.Range("E10").Value = "  -2.86%  "

$ws.Range("E11").Value = "  -0.77%  "

$ws.Range("E12").Value = "  +11.43%  "

$ws.Range("D13").Value = "'0.348"
$ws.Range("E13").Value = "  +2.11%  "

$ws.Range("D14").Value = "3.036.62"
$ws.Range("E14").Value = "  -2.81%  "

$ws.Range("D15").Value = "59.145.45"
$ws.Range("E15").Value = "  -2.11%  "

$ws.Range("D16").Value = "'23.11"
$ws.Range("E16").Value = "  +6.06%  "

$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").Value = "2.593.44"
$ws.Range("E18").Value = "  -2.17%  "

$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("D20").Value = "'336.83"
$ws.Range("E20").Value = "  -1.91%  "

$ws.Range("D21").Value = "'10.37"
$ws.Range("E21").Value = "  -0.39%  "

$ws.Range("D22").Value = "'6.39"
$ws.Range("E22").Value = "  -0.44%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").Value = "'64.06"
$ws.Range("E24").Value = "  -4.12%  "

$ws.Range("D25").Value = "'0.464"
$ws.Range("E25").Value = "  +4.84%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.30%  "

$ws.Range("E27").Value = "  -2.84%  "

$ws.Range("E28").Value = "  -0.23%  "

$ws.Range("D29").Value = "0.0₃0774"
$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("D31").Value = "'161.35"
$ws.Range("E31").Value = "  +3.35%  "

$ws.Range("E32").Value = "  -2.71%  "

$ws.Range("D33").Value = "'6.09"
$ws.Range("E33").Value = "  -0.65%  "

$ws.Range("E34").Value = "  -1.24%  "

$ws.Range("D35").Value = "'4.02"
$ws.Range("E35").Value = "  -1.63%  "

$ws.Range("E36").Value = "  -0.49%  "

$ws.Range("D37").Value = "'0.872"
$ws.Range("E37").Value = "  -3.16%  "

$ws.Range("D38").Value = "'0.876"
$ws.Range("E38").Value = "  -3.74%  "

$ws.Range("D39").Value = "'37.43"
$ws.Range("E39").Value = "  -0.19%  "

$ws.Range("E40").Value = "  -1.56%  "

$ws.Range("D41").Value = "'294.37"
$ws.Range("E41").Value = "  -3.39%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("E43").Value = "  +0.16%  "

$ws.Range("D44").Value = "'131.91"
$ws.Range("E44").Value = "  +5.10%  "

$ws.Range("E45").Value = "  -0.61%  "

$ws.Range("E46").Value = "  -1.43%  "

$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "'10.65"
$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0534"
$ws.Range("E48").Value = "  -2.48%  "

$ws.Range("D49").Value = "'19.02"
$ws.Range("E49").Value = "  -1.59%  "

$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("D51").Value = "'18.55"
$ws.Range("E51").Value = "  +0.48%  "

